$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Tool Usage Frequency" -> add Holder values in column E ---
$ws1 = $wb.Worksheets.Item("Tool Usage Frequency")

$holders = @{
    2  = "C50-10SM2"
    3  = "C50F3-10SF394-9"
    4  = "C50-50EM4"
    5  = "C50-75SF315-9"
    6  = "C50-16ER600"
    7  = "C50-16ER600"
    8  = "C50-38SF630-9"
    9  = "C50-75SF315-9"
    10 = "C50-32ER600"
    11 = "C50-50SF630-9"
    12 = "C50-50SF630-9"
    13 = "C50-32ER400"
    14 = "C50-10SM4"
    15 = "C50-75SF315-9"
    16 = "C50-50SF630-9"
    17 = "C50-16ER600"
    18 = "C50-16ER400"
    19 = "C50-16ER400"
    20 = "C50-16ER400"
    21 = "C50-62SF630-9"
    22 = "C50-32ER600"
    23 = "C50-32ER600"
    24 = "C50-32ER600"
    25 = "C50-16ER400"
    26 = "C50-25SF630-9"
    27 = "C50-32ER600"
    28 = ""
    29 = ""
    30 = ""
    31 = "C50-25SF630-9"
    32 = "C50-50SF630-9"
    33 = "C50-38SF630-9"
    34 = "C50-75SF315-9"
    35 = "C50-75SF630-9"
    36 = "C50-16ER600"
    37 = "C50-50SF630-9"
    38 = "C50-50SF630-9"
    39 = "C50-75SF315-9"
    40 = "C50-50SF630-9"
    41 = "C50-38SF630-9"
    42 = "C50-38SF315-9"
    43 = "C50-50SF630-9"
    44 = "C50-25SF630-9"
    45 = "C50-38SF315-9"
    46 = "C50-50SF315-9"
    47 = "C50-50SF630-9"
    48 = "C50-32ER600"
    49 = "C50-16ER400"
    50 = "C50-16ER400"
    51 = "C50-25SF630-9"
    52 = "C50-38SF630-9"
    53 = "C50-25SF630-9"
    54 = "C50-75SF630-9"
    55 = "C50-10EM4"
}

foreach ($row in 2..55) {
    $cell = $ws1.Cells.Item($row, 5)
    $cell.Value = $holders[$row]
    $cell.Style = "highlight"
}

# --- Sheet 3: "Programmer Stats" -> update counts / percentages ---
$ws3 = $wb.Worksheets.Item("Programmer Stats")

$ws3.Range("B2").Value = 38
$ws3.Range("C2").Value = 48.10126582278481

$ws3.Range("B3").Value = 40
$ws3.Range("C3").Value = 50.63291139240506

$ws3.Range("B4").Value = 1
$ws3.Range("C4").Value = 1.265822784810127
